# Refresh of the cryptos price table (Sheet1) -- updates Price (col D)
# and Volume(1h) (col E) figures, and for rows 37/38 also the Coin name
# and Link, matching the new ranking order (RenderToken now above
# LidoDAOToken).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price cells hold plain numeric-looking text (e.g. "68.40") in the
# workbook. Writing such a string straight into a General-formatted cell
# makes Excel reinterpret it as a number (dropping the trailing zero,
# e.g. 68.40 -> 68.4), so mark them as Text first to preserve the exact
# literal string, consistent with the rest of the column.
$forceTextCells = @("D5", "D6", "D9", "D10", "D11", "D12", "D13", "D14", "D16", "D18", "D20", "D22", "D23", "D24", "D26", "D28", "D29", "D30", "D32", "D33", "D35", "D37", "D38", "D39", "D44", "D48", "D49")
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Cells.Item(2, 4).Value = "43.396.21"
$ws.Cells.Item(2, 5).Value = "  +1.12%  "
$ws.Cells.Item(3, 4).Value = "2.372.39"
$ws.Cells.Item(3, 5).Value = "  +3.31%  "
$ws.Cells.Item(4, 5).Value = "  -0.02%  "
$ws.Cells.Item(5, 4).Value = "309.95"
$ws.Cells.Item(5, 5).Value = "  +0.16%  "
$ws.Cells.Item(6, 4).Value = "104.88"
$ws.Cells.Item(6, 5).Value = "  +5.17%  "
$ws.Cells.Item(7, 5).Value = "  -2.23%  "
$ws.Cells.Item(8, 5).Value = "  +0.01%  "
$ws.Cells.Item(9, 4).Value = "0.518"
$ws.Cells.Item(9, 5).Value = "  +0.56%  "
$ws.Cells.Item(10, 4).Value = "36.27"
$ws.Cells.Item(10, 5).Value = "  +0.87%  "
$ws.Cells.Item(11, 4).Value = "53.42"
$ws.Cells.Item(11, 5).Value = "  +2.65%  "
$ws.Cells.Item(12, 4).Value = "0.0813"
$ws.Cells.Item(12, 5).Value = "  -0.56%  "
$ws.Cells.Item(13, 4).Value = "0.113"
$ws.Cells.Item(13, 5).Value = "  -0.40%  "
$ws.Cells.Item(14, 4).Value = "7.01"
$ws.Cells.Item(14, 5).Value = "  -1.61%  "
$ws.Cells.Item(15, 4).Value = "2.741.31"
$ws.Cells.Item(16, 4).Value = "15.64"
$ws.Cells.Item(16, 5).Value = "  +5.39%  "
$ws.Cells.Item(17, 4).Value = "2.372.60"
$ws.Cells.Item(17, 5).Value = "  +2.80%  "
$ws.Cells.Item(18, 4).Value = "0.817"
$ws.Cells.Item(18, 5).Value = "  +2.31%  "
$ws.Cells.Item(19, 4).Value = "43.369.91"
$ws.Cells.Item(19, 5).Value = "  +1.15%  "
$ws.Cells.Item(20, 4).Value = "12.01"
$ws.Cells.Item(20, 5).Value = "  -3.32%  "
$ws.Cells.Item(21, 4).Value = "0.0₃0922"
$ws.Cells.Item(21, 5).Value = "  -0.17%  "
$ws.Cells.Item(22, 4).Value = "6.26"
$ws.Cells.Item(22, 5).Value = "  +3.44%  "
$ws.Cells.Item(23, 4).Value = "68.40"
$ws.Cells.Item(23, 5).Value = "  +0.50%  "
$ws.Cells.Item(24, 4).Value = "241.95"
$ws.Cells.Item(24, 5).Value = "  +1.11%  "
$ws.Cells.Item(25, 5).Value = "  +2.67%  "
$ws.Cells.Item(26, 4).Value = "2.62"
$ws.Cells.Item(26, 5).Value = "  +0.47%  "
$ws.Cells.Item(27, 5).Value = "  +0.24%  "
$ws.Cells.Item(28, 4).Value = "25.91"
$ws.Cells.Item(28, 5).Value = "  +7.74%  "
$ws.Cells.Item(29, 4).Value = "37.02"
$ws.Cells.Item(29, 5).Value = "  -3.90%  "
$ws.Cells.Item(30, 4).Value = "9.61"
$ws.Cells.Item(30, 5).Value = "  -0.13%  "
$ws.Cells.Item(31, 5).Value = "  +0.21%  "
$ws.Cells.Item(32, 4).Value = "162.05"
$ws.Cells.Item(32, 5).Value = "  -3.55%  "
$ws.Cells.Item(33, 4).Value = "5.28"
$ws.Cells.Item(33, 5).Value = "  -0.67%  "
$ws.Cells.Item(34, 5).Value = "  -0.02%  "
$ws.Cells.Item(35, 4).Value = "18.29"
$ws.Cells.Item(35, 5).Value = "  +3.50%  "
$ws.Cells.Item(36, 5).Value = "  +6.63%  "
$ws.Cells.Item(37, 2).Value = "RenderToken"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(37, 4).Value = "4.76"
$ws.Cells.Item(37, 5).Value = "  +12.67%  "
$ws.Cells.Item(38, 2).Value = "LidoDAOToken"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(38, 4).Value = "3.13"
$ws.Cells.Item(38, 5).Value = "  +0.35%  "
$ws.Cells.Item(39, 4).Value = "0.0743"
$ws.Cells.Item(39, 5).Value = "  +0.95%  "
$ws.Cells.Item(40, 5).Value = "  +6.42%  "
$ws.Cells.Item(41, 5).Value = "  +0.91%  "
$ws.Cells.Item(42, 5).Value = "  -1.26%  "
$ws.Cells.Item(43, 5).Value = "  +9.56%  "
$ws.Cells.Item(44, 4).Value = "19.95"
$ws.Cells.Item(44, 5).Value = "  +4.29%  "
$ws.Cells.Item(45, 4).Value = "2.005.83"
$ws.Cells.Item(45, 5).Value = "  +2.10%  "
$ws.Cells.Item(46, 5).Value = "  +0.93%  "
$ws.Cells.Item(47, 5).Value = "  +5.89%  "
$ws.Cells.Item(48, 4).Value = "10.47"
$ws.Cells.Item(48, 5).Value = "  +7.14%  "
$ws.Cells.Item(49, 4).Value = "58.34"
$ws.Cells.Item(49, 5).Value = "  +6.32%  "
$ws.Cells.Item(50, 5).Value = "  -0.01%  "
$ws.Cells.Item(51, 4).Value = "2.576.32"
$ws.Cells.Item(51, 5).Value = "  +2.06%  "
